# New crime data collected - weekly CompStat refresh (78th Precinct, week of 5/6/2024-5/12/2024)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Header text: volume/number and the reporting week dates (rich-text runs).
# Use Characters() to replace just the trailing numeric/date substrings.
# ---------------------------------------------------------------------------
$ws.Range("A8").Characters(21, 2).Text = "19"
# Replace the trailing date first (rightmost) so the earlier offset for the
# first date is unaffected by any length change from this edit.
$ws.Range("C9").Characters(47, 8).Text = "5/12/2024"
$ws.Range("C9").Characters(27, 9).Text = "5/6/2024"

# ---------------------------------------------------------------------------
# Helper reference cells whose style never changes across this edit - used
# as format+value donors via Range.Copy(destination) so the destination's
# type (numeric vs. shared-string placeholder) and style index switch
# together, matching how Excel performs such edits.
#   C14 -> s=14 text style, shared string "0"
#   E14 -> s=14 text style, shared string "***.*"
#   I14 -> s=15 numeric style
#   K14 -> s=16 percent style
# ---------------------------------------------------------------------------

# Row 14 (Murder): G14, H14 switch from numbers to text placeholders "0"/"***.*"
$ws.Range("C14").Copy($ws.Range("G14"))
$ws.Range("E14").Copy($ws.Range("H14"))

# Row 15 (Rape)
$ws.Range("N15").Value2 = -83.333333333333

# Row 16 (Robbery)
$ws.Range("C16").Value2 = 2
$ws.Range("D16").Value2 = 3
$ws.Range("E16").Value2 = -33.333333333333
$ws.Range("I16").Value2 = 33
$ws.Range("J16").Value2 = 48
$ws.Range("K16").Value2 = -31.25
$ws.Range("L16").Value2 = -8.333333333333
$ws.Range("M16").Value2 = -34
$ws.Range("N16").Value2 = -88.339222614841

# Row 17 (Fel. Assault)
$ws.Range("C17").Value2 = 4
$ws.Range("D17").Value2 = 1
$ws.Range("E17").Value2 = 300
$ws.Range("F17").Value2 = 10
$ws.Range("G17").Value2 = 7
$ws.Range("H17").Value2 = 42.857142857142
$ws.Range("I17").Value2 = 53
$ws.Range("J17").Value2 = 66
$ws.Range("K17").Value2 = -19.696969696969
$ws.Range("L17").Value2 = 32.5
$ws.Range("M17").Value2 = 430
$ws.Range("N17").Value2 = -28.378378378378

# Row 18 (Burglary)
$ws.Range("D18").Value2 = 2
$ws.Range("E18").Value2 = -50
$ws.Range("F18").Value2 = 6
$ws.Range("G18").Value2 = 10
$ws.Range("H18").Value2 = -40
$ws.Range("I18").Value2 = 43
$ws.Range("J18").Value2 = 79
$ws.Range("K18").Value2 = -45.569620253164
$ws.Range("L18").Value2 = 0
$ws.Range("M18").Value2 = -30.645161290322
$ws.Range("N18").Value2 = -83.773584905660

# Row 19 (Gr. Larceny)
$ws.Range("C19").Value2 = 9
$ws.Range("D19").Value2 = 11
$ws.Range("E19").Value2 = -18.181818181818
$ws.Range("F19").Value2 = 34
$ws.Range("G19").Value2 = 35
$ws.Range("H19").Value2 = -2.857142857142
$ws.Range("I19").Value2 = 156
$ws.Range("J19").Value2 = 199
$ws.Range("K19").Value2 = -21.608040201005
$ws.Range("L19").Value2 = -10.344827586206
$ws.Range("M19").Value2 = 4.697986577181
$ws.Range("N19").Value2 = 4

# Row 20 (G.L.A.): C20, D20, E20 switch from text placeholders to numbers
$ws.Range("I14").Copy($ws.Range("C20"))
$ws.Range("I14").Copy($ws.Range("D20"))
$ws.Range("K14").Copy($ws.Range("E20"))
$ws.Range("C20").Value2 = 1
$ws.Range("D20").Value2 = 4
$ws.Range("E20").Value2 = -75
$ws.Range("F20").Value2 = 12
$ws.Range("G20").Value2 = 7
$ws.Range("H20").Value2 = 71.428571428571
$ws.Range("I20").Value2 = 51
$ws.Range("J20").Value2 = 36
$ws.Range("K20").Value2 = 41.666666666666
$ws.Range("L20").Value2 = 88.888888888888
$ws.Range("M20").Value2 = 104
$ws.Range("N20").Value2 = -87.438423645320

# Row 21 (TOTAL)
$ws.Range("C21").Value2 = 17
$ws.Range("D21").Value2 = 21
$ws.Range("E21").Value2 = -19.047619047619
$ws.Range("F21").Value2 = 67
$ws.Range("G21").Value2 = 69
$ws.Range("H21").Value2 = -2.898550724637
$ws.Range("I21").Value2 = 340
$ws.Range("J21").Value2 = 431
$ws.Range("K21").Value2 = -21.113689095127
$ws.Range("L21").Value2 = 3.975535168195
$ws.Range("M21").Value2 = 14.093959731543
$ws.Range("N21").Value2 = -71.476510067114

# Row 22 (Transit): G22, H22 switch from numbers to text placeholders
$ws.Range("C14").Copy($ws.Range("G22"))
$ws.Range("E14").Copy($ws.Range("H22"))
$ws.Range("L22").Value2 = -10

# Row 23 (Housing): F23, G23, H23 switch from numbers to text placeholders
$ws.Range("C14").Copy($ws.Range("F23"))
$ws.Range("C14").Copy($ws.Range("G23"))
$ws.Range("E14").Copy($ws.Range("H23"))
$ws.Range("L23").Value2 = 40

# Row 24 (Petit Larceny)
$ws.Range("C24").Value2 = 22
$ws.Range("D24").Value2 = 41
$ws.Range("E24").Value2 = -46.341463414634
$ws.Range("F24").Value2 = 86
$ws.Range("G24").Value2 = 124
$ws.Range("H24").Value2 = -30.645161290322
$ws.Range("I24").Value2 = 385
$ws.Range("J24").Value2 = 583
$ws.Range("K24").Value2 = -33.962264150943
$ws.Range("L24").Value2 = -8.983451536643
$ws.Range("M24").Value2 = 35.563380281690

# Row 25 (Retail Theft)
$ws.Range("D25").Value2 = 28
$ws.Range("E25").Value2 = -50
$ws.Range("F25").Value2 = 57
$ws.Range("G25").Value2 = 89
$ws.Range("H25").Value2 = -35.955056179775
$ws.Range("I25").Value2 = 257
$ws.Range("J25").Value2 = 455
$ws.Range("K25").Value2 = -43.516483516483
$ws.Range("L25").Value2 = -24.633431085044

# Row 26 (Misd. Assault)
$ws.Range("C26").Value2 = 3
$ws.Range("D26").Value2 = 4
$ws.Range("E26").Value2 = -25
$ws.Range("F26").Value2 = 17
$ws.Range("G26").Value2 = 19
$ws.Range("H26").Value2 = -10.526315789473
$ws.Range("I26").Value2 = 74
$ws.Range("J26").Value2 = 108
$ws.Range("K26").Value2 = -31.481481481481
$ws.Range("L26").Value2 = -23.711340206185
$ws.Range("M26").Value2 = 7.246376811594

# Row 27 (UCR Rape*): F27, G27, H27 switch from numbers to text placeholders
$ws.Range("C14").Copy($ws.Range("F27"))
$ws.Range("C14").Copy($ws.Range("G27"))
$ws.Range("E14").Copy($ws.Range("H27"))

# Row 28 (Other Sex Crimes)
$ws.Range("C28").Value2 = 2
$ws.Range("I28").Value2 = 34
$ws.Range("K28").Value2 = 142.857142857143
$ws.Range("L28").Value2 = 88.888888888888

# Row 29 (Shooting Vic.): D29, E29, F29 switch from numbers to text placeholders
$ws.Range("C14").Copy($ws.Range("D29"))
$ws.Range("E14").Copy($ws.Range("E29"))
$ws.Range("C14").Copy($ws.Range("F29"))
$ws.Range("H29").Value2 = -100

# Row 30 (Shooting Inc.): D30, E30, F30 switch from numbers to text placeholders
$ws.Range("C14").Copy($ws.Range("D30"))
$ws.Range("E14").Copy($ws.Range("E30"))
$ws.Range("C14").Copy($ws.Range("F30"))
$ws.Range("H30").Value2 = -100

# Row 31 (Hate Crimes)
$ws.Range("G31").Value2 = 1
